# The document has two Pearson logos (in the footers) and two BTEC logos
# (in the headers). Their drawing-object names were swapped by mistake:
# the Pearson logos are named "image1.png" (should be "image2.png") and
# the BTEC logos are named "image2.jpg" (should be "image1.jpg"). Fix the
# InlineShape names accordingly.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTEC logo, currently "image2.jpg" -> "image1.jpg" ---
$hdrDefault = $sec.Headers.Item(1)
$hdrDefault.Range.InlineShapes.Item(1).Name = "image1.jpg"

$hdrFirstPage = $sec.Headers.Item(2)
$hdrFirstPage.Range.InlineShapes.Item(1).Name = "image1.jpg"

# --- Footers: Pearson logo, currently "image1.png" -> "image2.png" ---
# (Renaming a footer InlineShape directly off Footers(n).Range can leave
# the handle "stale"; selecting it first and going through the resulting
# Selection.InlineShapes collection avoids that.)
$ftrDefault = $sec.Footers.Item(1)
$ftrDefault.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

$ftrFirstPage = $sec.Footers.Item(2)
$ftrFirstPage.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"
